# Auto-generated edit script applying scheduled market-price refresh to Tonberry_Profits workbook.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 3329
$ws.Range("I52").Value = 3329
$ws.Range("K52").Value = 9987
$ws.Range("M52").Value = -9827
$ws.Range("H64").Value = 4578.3335
$ws.Range("I64").Value = 3700
$ws.Range("K64").Value = 3700
$ws.Range("M64").Value = -3452
$ws.Range("H67").Value = 4578.3335
$ws.Range("I67").Value = 3700
$ws.Range("K67").Value = 3700
$ws.Range("M67").Value = -2842
$ws.Range("H74").Value = 4225
$ws.Range("I74").Value = 3633.3333
$ws.Range("K74").Value = 3633.3333
$ws.Range("M74").Value = -2697.3333
$ws.Range("H77").Value = 4225
$ws.Range("I77").Value = 3633.3333
$ws.Range("K77").Value = 18166.6665
$ws.Range("M77").Value = -13486.6665
$ws.Range("H86").Value = 308755.75
$ws.Range("I86").Value = 616566.5
$ws.Range("J86").Value = 945
$ws.Range("K86").Value = 616566.5
$ws.Range("L86").Value = 945
$ws.Range("M86").Value = -615443.5
$ws.Range("N86").Value = -3191
$ws.Range("H89").Value = 308755.75
$ws.Range("I89").Value = 616566.5
$ws.Range("J89").Value = 945
$ws.Range("K89").Value = 3082832.5
$ws.Range("L89").Value = 4725
$ws.Range("M89").Value = -3077216.5
$ws.Range("N89").Value = -15957
$ws.Range("H112").Value = 5034.143
$ws.Range("J112").Value = 5689.8335
$ws.Range("L112").Value = 17069.5005
$ws.Range("N112").Value = -19285.5005
$ws.Range("H113").Value = 31857.143
$ws.Range("J113").Value = 3333.3333
$ws.Range("L113").Value = 3333.3333
$ws.Range("N113").Value = -9841.3333
$ws.Range("H116").Value = 12000.286
$ws.Range("I116").Value = 27750
$ws.Range("J116").Value = 5700.4
$ws.Range("K116").Value = 27750
$ws.Range("L116").Value = 5700.4
$ws.Range("M116").Value = -24308
$ws.Range("N116").Value = -12584.4
$ws.Range("H132").Value = 634.6667
$ws.Range("I132").Value = 573.87695
$ws.Range("K132").Value = 1721.63085
$ws.Range("M132").Value = 808.36915
$ws.Range("H138").Value = 2443.1667
$ws.Range("I138").Value = 2502
$ws.Range("J138").Value = 2373.6365
$ws.Range("K138").Value = 7506
$ws.Range("L138").Value = 7120.9095
$ws.Range("M138").Value = -2366
$ws.Range("N138").Value = -17400.9095
$ws.Range("H141").Value = 6317.3335
$ws.Range("I141").Value = 4000
$ws.Range("J141").Value = 6780.8
$ws.Range("K141").Value = 12000
$ws.Range("L141").Value = 20342.4
$ws.Range("M141").Value = -6820
$ws.Range("N141").Value = -30702.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1284.931
$ws.Range("I45").Value = 1097.1666
$ws.Range("K45").Value = 1097.1666
$ws.Range("M45").Value = -720.1666
$ws.Range("H122").Value = 1013.0769
$ws.Range("I122").Value = 1119.7222
$ws.Range("J122").Value = 773.125
$ws.Range("K122").Value = 3359.1666
$ws.Range("L122").Value = 2319.375
$ws.Range("M122").Value = -909.1665999999996
$ws.Range("N122").Value = -7219.375
$ws.Range("H132").Value = 2944.818
$ws.Range("I132").Value = 1832.8334
$ws.Range("K132").Value = 5498.5002
$ws.Range("M132").Value = -2968.5002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 20930.625
$ws.Range("I134").Value = 21935
$ws.Range("J134").Value = 13900
$ws.Range("K134").Value = 65805
$ws.Range("L134").Value = 41700
$ws.Range("M134").Value = -63270
$ws.Range("N134").Value = -46770

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1507
$ws.Range("I16").Value = 959.2
$ws.Range("K16").Value = 959.2
$ws.Range("M16").Value = -672.2
$ws.Range("H28").Value = 100000
$ws.Range("J28").Value = 100000
$ws.Range("L28").Value = 100000
$ws.Range("N28").Value = -100490
$ws.Range("H31").Value = 3092.8667
$ws.Range("I31").Value = 2953.3845
$ws.Range("K31").Value = 2953.3845
$ws.Range("M31").Value = -2658.3845
$ws.Range("H34").Value = 3092.8667
$ws.Range("I34").Value = 2953.3845
$ws.Range("K34").Value = 2953.3845
$ws.Range("M34").Value = -2751.3845
$ws.Range("H58").Value = 1116602
$ws.Range("I58").Value = 1673549.8
$ws.Range("K58").Value = 1673549.8
$ws.Range("M58").Value = -1673346.8
$ws.Range("H99").Value = 1927.091
$ws.Range("I99").Value = 1885.4286
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1885.4286
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -387.4286
$ws.Range("N99").Value = -4996
$ws.Range("H113").Value = 1507
$ws.Range("I113").Value = 959.2
$ws.Range("K113").Value = 959.2
$ws.Range("M113").Value = 1210.8
$ws.Range("H126").Value = 1927.091
$ws.Range("I126").Value = 1885.4286
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5656.2858
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3186.2858
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 1807.3948
$ws.Range("I132").Value = 995.24
$ws.Range("J132").Value = 3369.2307
$ws.Range("K132").Value = 2985.72
$ws.Range("L132").Value = 10107.6921
$ws.Range("M132").Value = -455.7200000000003
$ws.Range("N132").Value = -15167.6921
$ws.Range("H134").Value = 1652.6316
$ws.Range("I134").Value = 1652.6316
$ws.Range("K134").Value = 4957.8948
$ws.Range("M134").Value = -2422.8948
$ws.Range("H136").Value = 1116602
$ws.Range("I136").Value = 1673549.8
$ws.Range("K136").Value = 5020649.4
$ws.Range("M136").Value = -5018099.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 12001.571
$ws.Range("I87").Value = 6802.2
$ws.Range("K87").Value = 20406.6
$ws.Range("M87").Value = -19158.6
$ws.Range("H90").Value = 12001.571
$ws.Range("I90").Value = 6802.2
$ws.Range("K90").Value = 61219.8
$ws.Range("M90").Value = -54979.8
$ws.Range("H113").Value = 8454.308000000001
$ws.Range("J113").Value = 809.5454999999999
$ws.Range("L113").Value = 2428.6365
$ws.Range("N113").Value = -6768.6365
$ws.Range("H139").Value = 14460.875
$ws.Range("I139").Value = 18447.834
$ws.Range("J139").Value = 2500
$ws.Range("K139").Value = 55343.50199999999
$ws.Range("L139").Value = 7500
$ws.Range("M139").Value = -50203.50199999999
$ws.Range("N139").Value = -17780
$ws.Range("H140").Value = 2110.074
$ws.Range("I140").Value = 1726.5
$ws.Range("K140").Value = 5179.5
$ws.Range("M140").Value = 0.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2154.64
$ws.Range("I102").Value = 2212.182
$ws.Range("J102").Value = 1732.6666
$ws.Range("K102").Value = 2212.182
$ws.Range("L102").Value = 1732.6666
$ws.Range("M102").Value = -590.1819999999998
$ws.Range("N102").Value = -4976.6666
$ws.Range("H132").Value = 3499783
$ws.Range("I132").Value = 9617104
$ws.Range("K132").Value = 28851312
$ws.Range("M132").Value = -28848782

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2333.2666
$ws.Range("I22").Value = 2733.1667
$ws.Range("J22").Value = 2066.6667
$ws.Range("K22").Value = 2733.1667
$ws.Range("L22").Value = 2066.6667
$ws.Range("M22").Value = -2438.1667
$ws.Range("N22").Value = -2656.6667
$ws.Range("H27").Value = 2333.2666
$ws.Range("I27").Value = 2733.1667
$ws.Range("J27").Value = 2066.6667
$ws.Range("K27").Value = 2733.1667
$ws.Range("L27").Value = 2066.6667
$ws.Range("M27").Value = -2626.1667
$ws.Range("N27").Value = -2280.6667
$ws.Range("H136").Value = 5460.875
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = $null

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4525.6665
$ws.Range("I81").Value = 2788.5
$ws.Range("J81").Value = 8000
$ws.Range("K81").Value = 5577
$ws.Range("L81").Value = 16000
$ws.Range("M81").Value = -4516
$ws.Range("N81").Value = -18122
$ws.Range("H84").Value = 4525.6665
$ws.Range("I84").Value = 2788.5
$ws.Range("J84").Value = 8000
$ws.Range("K84").Value = 27885
$ws.Range("L84").Value = 80000
$ws.Range("M84").Value = -22581
$ws.Range("N84").Value = -90608
$ws.Range("H116").Value = 68000
$ws.Range("J116").Value = 68000
$ws.Range("L116").Value = 68000
$ws.Range("N116").Value = -77178
$ws.Range("H136").Value = 12921435
$ws.Range("I136").Value = 16836274
$ws.Range("K136").Value = 50508822
$ws.Range("M136").Value = -50506272
